$wb = $excel.ActiveWorkbook

# "Jack" sheet: delete the first data row (row 2, year 2020). Every
# subsequent year shifts up one row, effectively bumping all years by one
# (2021 -> row2, 2022 -> row3, ... 2055 -> row36) and the sheet shrinks by a
# row (A1:I37 -> A1:I36).
$wsJack = $wb.Worksheets.Item("Jack")
$wsJack.Activate()
$wsJack.Rows("2").Delete() | Out-Null
$wsJack.Range("A2:A36").Select() | Out-Null

# "Jill" sheet: same row-shift operation (A1:I40 -> A1:I39), ending up as the
# active sheet/tab with the selection left on E12.
$wsJill = $wb.Worksheets.Item("Jill")
$wsJill.Activate()
$wsJill.Rows("2").Delete() | Out-Null
$wsJill.Range("E12").Select() | Out-Null
